$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row to add at the bottom of the localization table
$newRow = 54

# Set the text values for the new row first
$ws.Cells.Item($newRow, 1).Value = "lang_select_place"
$ws.Cells.Item($newRow, 2).Value = "Chọn Địa Chỉ Này"
$ws.Cells.Item($newRow, 3).Value = "Select This Address"

# Column A on this sheet uses a distinct "key" style (Consolas, vertical
# centered). Copy that formatting from the row above so we reuse the
# existing font/style entries instead of creating new duplicate ones.
$ws.Cells.Item($newRow - 1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

# Columns B and C reuse the plain style already used throughout the sheet.
$ws.Cells.Item($newRow - 1, 2).Copy()
$ws.Cells.Item($newRow, 2).PasteSpecial(-4122)

$ws.Cells.Item($newRow - 1, 3).Copy()
$ws.Cells.Item($newRow, 3).PasteSpecial(-4122)

# Give the new key cell its own accent color (matches a VS Code string
# token color) distinguishing it as a newly added localization key.
$ws.Cells.Item($newRow, 1).Font.Color = 7901646

# Keep the current selection on the newly added cell, like Excel would
# after typing in the last row.
$ws.Cells.Item($newRow, 3).Select()
